# إضافة حدث جديد في Card23
# Fills in the blank tracking cells of the last existing event row (17)
# with "nan" placeholders (matching the sheet's convention for "no data"),
# and appends a brand-new service event as row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Helper: force a cell to be stored as TEXT (not auto-converted to a
# number) even when the value looks numeric (e.g. "23", "641.6"), then
# drop the temporary "Text" number format again so the cell ends up with
# the workbook's normal (unstyled) look.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Row 17: the previously-empty tracking columns now read "nan" ---
$row17Cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "M")
foreach ($col in $row17Cols) {
    $ws.Range($col + "17").Value = "nan"
}

# --- Row 18: brand-new service event for Card23 ---
Set-TextValue $ws.Range("A18") "23"

# B18..K18 stay blank (same as the other "no measurement" columns),
# but keep them as present, unstyled cells like the rest of the sheet.
$blankCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
foreach ($col in $blankCols) {
    $cell = $ws.Range($col + "18")
    $cell.NumberFormat = "@"
    $cell.ClearFormats()
}

$ws.Range("L18").Value = "15\5\2025"
Set-TextValue $ws.Range("M18") "641.6"
$ws.Range("N18").Value = "تم عمل صيانه وتغيير الجرائد الاماميه ومعايره المكنه (1_2_4_5_7_8)"
$ws.Range("O18").Value = "الخبير"
